$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the missing "Status" value for row 6
$ws.Range("F6").Value = "PASS"

# Add a new row of login data (row 7)
$ws.Range("C7").Value = "sarzigptgnrrjrf@gmail.com"
$ws.Range("D7").Value = "wcsjcRVJTP5"
$ws.Range("E7").Value = "pass"
$ws.Range("F7").Value = "PASS"
